# faturamento_diario.xlsx update:
#  1. Insert a new daily record (day 24, July/2025) as a new row 25,
#     pushing the existing rows 25-115 down to 26-116.
#  2. Correct the sales total for day 23 (July/2025), currently in row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25 (shifts everything below it down by one)
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row with the new daily record
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = 15029.97
$ws.Cells.Item(25, 3).Value = 7
$ws.Cells.Item(25, 4).Value = 2025
$ws.Cells.Item(25, 5).Value = "07/2025"

# Update the corrected total for day 23 (row 24)
$ws.Cells.Item(24, 2).Value = 14800.32
